# Release 4.2.1: Update Excel files, fixes
# - Replace the STIX ID (column B) UUIDs for every tactic row with new ones.
# - Fix a few typos/wording issues in the tactic descriptions (column D).
# - Bump the "last modified" date (shared by created/last modified columns F & G)
#   from "11 March 2024" to "29 October 2024".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tactics")

# --- Update STIX IDs in column B (rows 2-15) ---
$ws.Range("B2").Value  = "x-mitre-tactic--76d8d1bc-98fb-406a-a88d-a70649ad9365"
$ws.Range("B3").Value  = "x-mitre-tactic--b7c3883a-1ae7-47cd-bd53-024c79ac833a"
$ws.Range("B4").Value  = "x-mitre-tactic--35a4b685-8a9f-4f66-bef0-92c8a5fd8411"
$ws.Range("B5").Value  = "x-mitre-tactic--b3f5fda4-a31a-4a74-9588-995b143d3436"
$ws.Range("B6").Value  = "x-mitre-tactic--b777d4af-fbe8-4e44-9917-c5c5542a7147"
$ws.Range("B7").Value  = "x-mitre-tactic--55e80a13-be94-43ba-aa57-e6c5cb822864"
$ws.Range("B8").Value  = "x-mitre-tactic--4386734e-fa57-4700-95b8-e76b8cab8ab3"
$ws.Range("B9").Value  = "x-mitre-tactic--14ef1e2d-f902-450a-a7e8-b032e318bb86"
$ws.Range("B10").Value = "x-mitre-tactic--528bc2da-8855-44c0-9e45-ab92c179bf6f"
$ws.Range("B11").Value = "x-mitre-tactic--2666e534-bf47-4656-9404-e90afe41f4a7"
$ws.Range("B12").Value = "x-mitre-tactic--4a088a36-6786-4486-a3ea-3c3576d61daa"
$ws.Range("B13").Value = "x-mitre-tactic--ef06a48d-1ccb-42c6-b5dd-2770a58f02d8"
$ws.Range("B14").Value = "x-mitre-tactic--9ceaa8fe-57f1-4923-a4a5-121b5111139c"
$ws.Range("B15").Value = "x-mitre-tactic--1fa998a1-d720-488c-b156-b569fe4e6308"

# --- Fix wording in descriptions (column D) ---
# Apply targeted text substitutions on top of the existing cell contents so that
# unrelated formatting (e.g. blank lines already present between the first
# sentence and the following paragraph) is preserved untouched.

# Row 6: Execution - turn "Remote System Discovery" into a Markdown link
$cell = $ws.Range("D6")
$text = $cell.Value()
$text = $text.Replace( `
    "does Remote System Discovery.", `
    "does [Remote System Discovery](https://attack.mitre.org/techniques/T1018/).")
$cell.Value = $text

# Row 10: ML Attack Staging - "manor" -> "manner"
$cell = $ws.Range("D10")
$text = $cell.Value()
$text = $text.Replace( `
    "offline manor and are thus difficult to mitigate.", `
    "offline manner and are thus difficult to mitigate.")
$cell.Value = $text

# Row 14: Reconnaissance - "organizations" -> "organizations'"
$cell = $ws.Range("D14")
$text = $cell.Value()
$text = $text.Replace( `
    "victim organizations machine learning capabilities", `
    "victim organizations' machine learning capabilities")
$cell.Value = $text

# Row 15: Resource Development - turn "ML Attack Staging" into a relative link
$cell = $ws.Range("D15")
$text = $cell.Value()
$text = $text.Replace( `
    "lifecycle, such as ML Attack Staging.", `
    "lifecycle, such as [ML Attack Staging](/tactics/AML.TA0001).")
$cell.Value = $text

# --- Bump the "last modified" / "created" date shared string for all rows (F2:G15) ---
for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 6).Value = "29 October 2024"
    $ws.Cells.Item($r, 7).Value = "29 October 2024"
}
